# Automatische test-sync: 2025-08-19 20:25:50
# Append a new log entry to the "Logs" sheet and refresh the dependent
# conditional formatting ranges plus the "Dashboard" summary count.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append the new row of log data (row 17) ---------------------------
$newRow = 17

$logs.Range("A" + $newRow).Value = "Interne taak"
$logs.Range("B" + $newRow).Value = "kwaliteit@testbedrijf123.nl"
$logs.Range("D" + $newRow).Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F" + $newRow).Value = "2025-08-19 20:25:36"
$logs.Range("G" + $newRow).Value = "Nee"
$logs.Range("H" + $newRow).Value = "Ja"
$logs.Range("I" + $newRow).Value = "Nee"
$logs.Range("J" + $newRow).Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row ---
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range($col + "2:" + $col + "16")
    $newRange = $logs.Range($col + "2:" + $col + $newRow)
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Update the Dashboard summary count ---------------------------------
$dashboard.Range("B2").Value = 16
